$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$years = @(2021, 2022, 2023, 2024, 2025, 2026, 2027, 2028, 2029, 2030)
$rice  = @(73336.02, 73966.17, 74596.26, 75226.32, 75856.32, 76486.55, 77014.89, 77543.28, 78071.96, 78600.27)
$wheat = @(54750.39, 55238.76, 55727.17, 56215.43, 56703.76, 57192.15, 57601.68, 58011.2, 58420.9, 58830.37)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $rice[$i]
    $ws.Cells.Item($row, 3).Value = $wheat[$i]
}
